$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 160.57143
$ws.Range("I2").Value = 89.59999999999999
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 89.59999999999999
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 23.40000000000001
$ws.Range("N2").Value = -426

$ws.Range("H70").Value = 14599.889
$ws.Range("I70").Value = 4850
$ws.Range("J70").Value = 22399.8
$ws.Range("K70").Value = 14550
$ws.Range("L70").Value = 67199.39999999999
$ws.Range("M70").Value = -14280
$ws.Range("N70").Value = -67739.39999999999

$ws.Range("H73").Value = 14599.889
$ws.Range("I73").Value = 4850
$ws.Range("J73").Value = 22399.8
$ws.Range("K73").Value = 14550
$ws.Range("L73").Value = 67199.39999999999
$ws.Range("M73").Value = -13614
$ws.Range("N73").Value = -69071.39999999999

$ws.Range("H86").Value = 3744.3333
$ws.Range("J86").Value = 3493.2
$ws.Range("L86").Value = 3493.2
$ws.Range("N86").Value = -5739.2

$ws.Range("H89").Value = 3744.3333
$ws.Range("J89").Value = 3493.2
$ws.Range("L89").Value = 17466
$ws.Range("N89").Value = -28698

$ws.Range("H99").Value = 671
$ws.Range("J99").Value = 999.3333
$ws.Range("L99").Value = 2997.9999
$ws.Range("N99").Value = -5993.9999

$ws.Range("H103").Value = 713.5
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 618
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 1854
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -3026

$ws.Range("H135").Value = 1438
$ws.Range("I135").Value = 1319.3334
$ws.Range("K135").Value = 11874.0006
$ws.Range("M135").Value = -9339.000599999999

$ws.Range("H137").Value = 501607.34
$ws.Range("I137").Value = 1001159
$ws.Range("J137").Value = 2055.7
$ws.Range("K137").Value = 3003477
$ws.Range("L137").Value = 6167.099999999999
$ws.Range("M137").Value = -3000927
$ws.Range("N137").Value = -11267.1

$ws.Range("H141").Value = 4961.4585
$ws.Range("I141").Value = 4023.158
$ws.Range("J141").Value = 8527
$ws.Range("K141").Value = 12069.474
$ws.Range("L141").Value = 25581
$ws.Range("M141").Value = -6889.474
$ws.Range("N141").Value = -35941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 10185.917
$ws.Range("I88").Value = 1690
$ws.Range("J88").Value = 13017.889
$ws.Range("K88").Value = 1690
$ws.Range("L88").Value = 13017.889
$ws.Range("M88").Value = -1284
$ws.Range("N88").Value = -13829.889

$ws.Range("H91").Value = 10185.917
$ws.Range("I91").Value = 1690
$ws.Range("J91").Value = 13017.889
$ws.Range("K91").Value = 1690
$ws.Range("L91").Value = 13017.889
$ws.Range("M91").Value = -286
$ws.Range("N91").Value = -15825.889

$ws.Range("H97").Value = 3241.878
$ws.Range("I97").Value = 3762.2258
$ws.Range("J97").Value = 1628.8
$ws.Range("K97").Value = 3762.2258
$ws.Range("L97").Value = 1628.8
$ws.Range("M97").Value = -3266.2258
$ws.Range("N97").Value = -2620.8

$ws.Range("H102").Value = 1395.2858
$ws.Range("I102").Value = 1383.6316
$ws.Range("K102").Value = 1383.6316
$ws.Range("M102").Value = 238.3684000000001

$ws.Range("H122").Value = 3043.36
$ws.Range("I122").Value = 2694.2778
$ws.Range("J122").Value = 3941
$ws.Range("K122").Value = 8082.8334
$ws.Range("L122").Value = 11823
$ws.Range("M122").Value = -5632.8334
$ws.Range("N122").Value = -16723

$ws.Range("H132").Value = 1993.5454
$ws.Range("I132").Value = 1449
$ws.Range("K132").Value = 4347
$ws.Range("M132").Value = -1817

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 15000
$ws.Range("J9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("N9").Value = -15336

$ws.Range("H22").Value = 1818907.9
$ws.Range("I22").Value = 316.33334
$ws.Range("K22").Value = 316.33334
$ws.Range("M22").Value = -143.33334

$ws.Range("H26").Value = 39181
$ws.Range("I26").Value = 39181
$ws.Range("K26").Value = 39181
$ws.Range("M26").Value = -38889

$ws.Range("H134").Value = 8100.5293
$ws.Range("I134").Value = 7405.7085
$ws.Range("K134").Value = 22217.1255
$ws.Range("M134").Value = -19682.1255

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1706.6428
$ws.Range("J94").Value = 1821.25
$ws.Range("L94").Value = 1821.25
$ws.Range("N94").Value = -2723.25

$ws.Range("H105").Value = 1144.3636
$ws.Range("I105").Value = 1069.1
$ws.Range("K105").Value = 1069.1
$ws.Range("M105").Value = 677.9000000000001

$ws.Range("H134").Value = 4496.9023
$ws.Range("I134").Value = 3780.8572
$ws.Range("K134").Value = 11342.5716
$ws.Range("M134").Value = -8807.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1688617.9
$ws.Range("I11").Value = 3658322
$ws.Range("K11").Value = 10974966
$ws.Range("M11").Value = -10974826

$ws.Range("H26").Value = 750.6
$ws.Range("I26").Value = 251.5
$ws.Range("J26").Value = 1083.3334
$ws.Range("K26").Value = 754.5
$ws.Range("L26").Value = 3250.0002
$ws.Range("M26").Value = -466.5
$ws.Range("N26").Value = -3826.0002

$ws.Range("H58").Value = 10003
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 10003
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 30009
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -30265

$ws.Range("H117").Value = 4798
$ws.Range("J117").Value = 6330
$ws.Range("L117").Value = 18990
$ws.Range("N117").Value = -25874

$ws.Range("H131").Value = 4731625.5
$ws.Range("I131").Value = 10101916
$ws.Range("J131").Value = 4019900
$ws.Range("K131").Value = 30305748
$ws.Range("L131").Value = 12059700
$ws.Range("M131").Value = -30300708
$ws.Range("N131").Value = -12069780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2182.7778
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 2377.8572
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 2377.8572
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -4373.8572

$ws.Range("H83").Value = 2182.7778
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 2377.8572
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 11889.286
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -21873.286

$ws.Range("H102").Value = 9975
$ws.Range("I102").Value = 11633.333
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 11633.333
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -10011.333
$ws.Range("N102").Value = -8244

$ws.Range("H113").Value = 3049.6667
$ws.Range("I113").Value = 3474.5
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 3474.5
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -1304.5
$ws.Range("N113").Value = -6540

$ws.Range("H132").Value = 2705.85
$ws.Range("I132").Value = 2672.7222
$ws.Range("J132").Value = 3004
$ws.Range("K132").Value = 8018.1666
$ws.Range("L132").Value = 9012
$ws.Range("M132").Value = -5488.1666
$ws.Range("N132").Value = -14072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1972.625
$ws.Range("J46").Value = 2710.889
$ws.Range("L46").Value = 2710.889
$ws.Range("N46").Value = -3086.889

$ws.Range("H61").Value = 1313.8572
$ws.Range("I61").Value = 1338
$ws.Range("K61").Value = 1338
$ws.Range("M61").Value = -1136

$ws.Range("H82").Value = 3364.3
$ws.Range("I82").Value = 1644.3846
$ws.Range("K82").Value = 1644.3846
$ws.Range("M82").Value = -1283.3846

$ws.Range("H85").Value = 3364.3
$ws.Range("I85").Value = 1644.3846
$ws.Range("K85").Value = 1644.3846
$ws.Range("M85").Value = -396.3846000000001

$ws.Range("H100").Value = 1748
$ws.Range("I100").Value = 1598
$ws.Range("K100").Value = 1598
$ws.Range("M100").Value = -1057

$ws.Range("H113").Value = 1313.8572
$ws.Range("I113").Value = 1338
$ws.Range("K113").Value = 1338
$ws.Range("M113").Value = 832

$ws.Range("H122").Value = 5042.483
$ws.Range("I122").Value = 2912.5454
$ws.Range("J122").Value = 6344.1113
$ws.Range("K122").Value = 8737.636200000001
$ws.Range("L122").Value = 19032.3339
$ws.Range("M122").Value = -6287.636200000001
$ws.Range("N122").Value = -23932.3339

$ws.Range("H132").Value = 3365.5
$ws.Range("I132").Value = 2557
$ws.Range("K132").Value = 7671
$ws.Range("M132").Value = -5141

$ws.Range("H136").Value = 2523.9788
$ws.Range("I136").Value = 1197.8096
$ws.Range("J136").Value = 3595.1155
$ws.Range("K136").Value = 3593.4288
$ws.Range("L136").Value = 10785.3465
$ws.Range("M136").Value = -1043.4288
$ws.Range("N136").Value = -15885.3465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1164.2858
$ws.Range("I100").Value = 553.5
$ws.Range("K100").Value = 1107
$ws.Range("M100").Value = -566

